# 17 Mayis 2020 verileri eklendi
# Adds the 17 May 2020 COVID-19 data row (row 67) to the "data" worksheet,
# expands Table3 to include the new row, and updates the active selection
# to match the recorded state after entry (E66, since the new row below it
# becomes the insertion point for the next day's data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: date=17 May 2020 (serial 43968), test, case, death, recovered
$ws.Range("A67").Value = 43968
$ws.Range("B67").Value = 35369
$ws.Range("C67").Value = 1368
$ws.Range("D67").Value = 44
$ws.Range("E67").Value = 1825

# Expand the worksheet table ("Table3") so its range / autofilter covers the
# newly added row.
$tbl = $ws.ListObjects.Item("Table3")
$tbl.Resize($ws.Range("A1:E67"))

# Match the saved selection state (cell below the previous last row).
$ws.Range("E66").Select()
